# "added clock face intialization"
# Populate the small helper table in J18:N29 that computes the clock-face
# tick-mark coordinates (M/N) from a radius (K18) and an angle column (L),
# re-using the existing number/alignment formatting already present on the
# sheet (columns A and G) so no redundant styles are introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J18: an empty cell that merely carries an integer number format (no
# alignment) - this is a brand new style, so it becomes the new cellXfs entry.
$ws.Range("J18").NumberFormat = "0"

# K18: radius used by every row of the table below.
$ws.Range("K18").Value = 45

# L18:L29: the angle, in degrees, for each of the 12 clock positions.
# Formatting copied from A18 (center/center alignment, no number format).
$angles = @(60, 30, 0, -30, -60, -90, -120, -150, -180, -210, -240, -270)
for ($i = 0; $i -lt $angles.Length; $i++) {
    $row = 18 + $i
    $ws.Range("L$row").Value = $angles[$i]
}
$ws.Range("A18").Copy()
$ws.Range("L18:L29").PasteSpecial(-4122)

# M18:M29 / N18:N29: the x/y pixel coordinates for each tick mark, driven by
# the radius in K18 and the angle in the same row's L cell. Formatting
# copied from G18 (center/center alignment, integer number format).
for ($row = 18; $row -le 29; $row++) {
    $ws.Range("M$row").Formula = "=`$K`$18*COS(RADIANS(`$L$row))+`$J`$3"
    $ws.Range("N$row").Formula = "=`$J`$4-`$K`$18*SIN(RADIANS(`$L$row))-4"
}
$ws.Range("G18").Copy()
$ws.Range("M18:N29").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# The author's selection ended on the last newly-entered pair of cells.
$ws.Range("M29:N29").Select() | Out-Null
